# Apply "Doing Updates for Financials" edit:
# Insert a new first data column (D) for the latest reporting period (2018-09-30, serial 43343)
# on the ACN worksheet, shifting the existing period columns one column to the right, and
# correct a handful of values whose originally-reported figures were revised compared with
# what the mechanical shift alone would produce.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new column before D; this shifts D:K -> E:L
# ------------------------------------------------------------------
$ws.Columns("D").Insert()

# ------------------------------------------------------------------
# 2. Give the new column D the same formatting as column E (which holds
#    what used to be column D), for every row that actually has data.
# ------------------------------------------------------------------
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Populate the new column D with the latest period's figures.
# ------------------------------------------------------------------
$dCol = @{
    7 = 43343
    8 = 41603400
    9 = 29160500
    10 = 12442900
    12 = "NA"
    13 = 0
    14 = 0
    15 = 0
    17 = 35762400
    18 = 5841000
    20 = -13400
    21 = 6754400
    22 = 19500
    23 = 5808100
    24 = 1415800
    25 = 0
    26 = 4392200
    27 = 4237600
    28 = 0
    29 = -177700
    30 = 0
    31 = 0
    32 = 13400
    33 = 4059900
    34 = 0
    35 = 4059900
    38 = 43343
    41 = 5015700
    42 = 3200
    43 = 7496400
    44 = 0
    45 = 1070300
    46 = 13585600
    47 = 238600
    48 = 1264000
    49 = 6070100
    50 = 0
    51 = 0
    52 = 3290800
    53 = 0
    54 = 24449100
    57 = 1348800
    58 = 5300
    59 = 8797600
    60 = 10151800
    61 = 19700
    62 = 3553100
    63 = 0
    64 = 0
    65 = 0
    66 = 14084300
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 7952400
    73 = 0
    74 = 0
    75 = 0
    76 = 10364800
    77 = 0
    80 = 43343
    81 = 4059900
    83 = 926800
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 6026700
    91 = -619200
    92 = 0
    93 = 0
    94 = -1249600
    96 = -1671100
    97 = 0
    98 = 0
    99 = 0
    100 = -3709000
    101 = -133600
    102 = 934500
}

foreach ($r in $dCol.Keys) {
    $ws.Cells.Item($r, 4).Value = $dCol[$r]
}

# ------------------------------------------------------------------
# 4. A handful of rows received genuinely revised figures (not just a
#    column shift) for the most recent periods; overwrite those cells
#    with the corrected values.
# ------------------------------------------------------------------
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = 0

$ws.Range("E41").Value = 4081300
$ws.Range("F41").Value = 4860100
$ws.Range("G41").Value = 8675600

$ws.Range("E45").Value = 1127700
$ws.Range("F45").Value = 890800
$ws.Range("G45").Value = 1536700

$ws.Range("E96").Value = -1498700
$ws.Range("F96").Value = -1372200
$ws.Range("G96").Value = -1276800

Write-Host "Applied financial update edit"
